# Refactor timetable generation to improve formatting and readability by
# adding line breaks in activity descriptions (Tomasz SKWERES Camp A timetable).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3: rename the class and switch "&" to "/" while keeping the room line break.
$ws.Range("D3").Value = "Cello Regulation / Maintance Workshop`n(Room 201)"

# Insert "Private" before "Lesson" in the one-on-one lesson cells, and
# normalise the double space before "Lesson" to a single space where present.
$ws.Range("C7").Value  = "Henry SUN Private Lesson with Tomasz SKWERES & pianist"
$ws.Range("E7").Value  = "Yauyau NG Private Lesson with Tomasz SKWERES & pianist"
$ws.Range("C11").Value = "Peter CHAN Private Lesson with Tomasz SKWERES & pianist"
$ws.Range("E11").Value = "Icelyn GE Private Lesson with Tomasz SKWERES & pianist"
$ws.Range("C19").Value = "Max TAM Private Lesson with Tomasz SKWERES & pianist"
$ws.Range("E19").Value = "Ivy YUE Private Lesson with Tomasz SKWERES & pianist"
